# Weekly update: insert a new observation row for Femacal de La Calera - Poroto verde
# This shifts the existing rows 215-266 down to 216-267 and populates the new row 215
# with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 215 (pushes rows 215:266 down to 216:267)
$ws.Rows(215).Insert()

# Populate the newly inserted row 215 with the new weekly record
$ws.Range("A215").Value = 3
$ws.Range("B215").Value = "Femacal de La Calera"
$ws.Range("C215").Value = "Coquimbo"
$ws.Range("D215").Value2 = 44543
$ws.Range("E215").Value = 5
$ws.Range("F215").Value = 100112031
$ws.Range("G215").Value = "Poroto verde"
$ws.Range("H215").Value = "Magnum"
$ws.Range("I215").Value = "Primera"
$ws.Range("J215").Value = 76
$ws.Range("K215").Value = 20000
$ws.Range("L215").Value = 21000
$ws.Range("M215").Value = 20500
$ws.Range("N215").Value = "`$/malla 25 kilos"
$ws.Range("O215").Value = "Provincia de Limarí"
$ws.Range("P215").Value = 820
$ws.Range("Q215").Value = 25
$ws.Range("R215").Value = "Hortaliza"
